$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 68
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = -290
$ws.Range("H6").Value = 1105.091
$ws.Range("I6").Value = 131.2
$ws.Range("J6").Value = 1916.6666
$ws.Range("K6").Value = 393.6
$ws.Range("L6").Value = 5749.9998
$ws.Range("M6").Value = -281.6
$ws.Range("N6").Value = -5973.9998
$ws.Range("H17").Value = 6758.154
$ws.Range("J17").Value = 1672.8422
$ws.Range("L17").Value = 5018.5266
$ws.Range("N17").Value = -5354.5266
$ws.Range("H31").Value = 1200
$ws.Range("I31").Value = 1200
$ws.Range("K31").Value = 3600
$ws.Range("M31").Value = -3370
$ws.Range("H107").Value = 727.75
$ws.Range("I107").Value = 755.7368
$ws.Range("J107").Value = 196
$ws.Range("K107").Value = 755.7368
$ws.Range("L107").Value = 196
$ws.Range("M107").Value = 1164.2632
$ws.Range("N107").Value = -4036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11906136
$ws.Range("I2").Value = 16667405
$ws.Range("J2").Value = 2962.6667
$ws.Range("K2").Value = 16667405
$ws.Range("L2").Value = 2962.6667
$ws.Range("M2").Value = -16667292
$ws.Range("N2").Value = -3188.6667
$ws.Range("H25").Value = 30885.625
$ws.Range("I25").Value = 1406.2
$ws.Range("K25").Value = 1406.2
$ws.Range("M25").Value = -1004.2
$ws.Range("H51").Value = 50023.5
$ws.Range("J51").Value = 50023.5
$ws.Range("L51").Value = 50023.5
$ws.Range("N51").Value = -51535.5
$ws.Range("H116").Value = 11906136
$ws.Range("I116").Value = 16667405
$ws.Range("J116").Value = 2962.6667
$ws.Range("K116").Value = 16667405
$ws.Range("L116").Value = 2962.6667
$ws.Range("M116").Value = -16665111
$ws.Range("N116").Value = -7550.6667
$ws.Range("H132").Value = 16395795
$ws.Range("I132").Value = 25001288
$ws.Range("K132").Value = 75003864
$ws.Range("M132").Value = -75001334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11906136
$ws.Range("I3").Value = 16667405
$ws.Range("J3").Value = 2962.6667
$ws.Range("K3").Value = 16667405
$ws.Range("L3").Value = 2962.6667
$ws.Range("M3").Value = -16667291
$ws.Range("N3").Value = -3190.6667
$ws.Range("H22").Value = 260
$ws.Range("I22").Value = 96
$ws.Range("J22").Value = 533.3333
$ws.Range("K22").Value = 96
$ws.Range("L22").Value = 533.3333
$ws.Range("M22").Value = 77
$ws.Range("N22").Value = -879.3333
$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178
$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 169.91667
$ws.Range("I7").Value = 212.71428
$ws.Range("J7").Value = 110
$ws.Range("K7").Value = 212.71428
$ws.Range("L7").Value = 110
$ws.Range("M7").Value = -99.71428
$ws.Range("N7").Value = -336
$ws.Range("H10").Value = 5057.45
$ws.Range("I10").Value = 746.2308
$ws.Range("J10").Value = 13064
$ws.Range("K10").Value = 746.2308
$ws.Range("L10").Value = 13064
$ws.Range("M10").Value = -607.2308
$ws.Range("N10").Value = -13342
$ws.Range("H16").Value = 1529.4
$ws.Range("I16").Value = 867
$ws.Range("J16").Value = 2191.8
$ws.Range("K16").Value = 867
$ws.Range("L16").Value = 2191.8
$ws.Range("M16").Value = -580
$ws.Range("N16").Value = -2765.8
$ws.Range("H21").Value = 20958.637
$ws.Range("H58").Value = 9617886
$ws.Range("I58").Value = 1514.2059
$ws.Range("J58").Value = 27782144
$ws.Range("K58").Value = 1514.2059
$ws.Range("L58").Value = 27782144
$ws.Range("M58").Value = -1311.2059
$ws.Range("N58").Value = -27782550
$ws.Range("H59").Value = 22751.334
$ws.Range("J59").Value = 22751.334
$ws.Range("L59").Value = 22751.334
$ws.Range("N59").Value = -25041.334
$ws.Range("H105").Value = 5900.5
$ws.Range("I105").Value = 6143.5713
$ws.Range("K105").Value = 6143.5713
$ws.Range("M105").Value = -4396.5713
$ws.Range("H113").Value = 1529.4
$ws.Range("I113").Value = 867
$ws.Range("J113").Value = 2191.8
$ws.Range("K113").Value = 867
$ws.Range("L113").Value = 2191.8
$ws.Range("M113").Value = 1303
$ws.Range("N113").Value = -6531.8
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("H136").Value = 9617886
$ws.Range("I136").Value = 1514.2059
$ws.Range("J136").Value = 27782144
$ws.Range("K136").Value = 4542.6177
$ws.Range("L136").Value = 83346432
$ws.Range("M136").Value = -1992.6177
$ws.Range("N136").Value = -83351532

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 21607.715
$ws.Range("I4").Value = 27364.363
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 82093.08900000001
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -81981.08900000001
$ws.Range("N4").Value = -1724
$ws.Range("H5").Value = 1237.55
$ws.Range("I5").Value = 764.7143
$ws.Range("K5").Value = 2294.1429
$ws.Range("M5").Value = -2182.1429
$ws.Range("H7").Value = 166.33333
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 299
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 897
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -1121
$ws.Range("H122").Value = 2038.2858
$ws.Range("I122").Value = 511
$ws.Range("J122").Value = 2886.7778
$ws.Range("K122").Value = 4599
$ws.Range("L122").Value = 25981.0002
$ws.Range("M122").Value = -2149
$ws.Range("N122").Value = -30881.0002
$ws.Range("H135").Value = 1237.55
$ws.Range("I135").Value = 764.7143
$ws.Range("K135").Value = 6882.428699999999
$ws.Range("M135").Value = -4347.428699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.263159
$ws.Range("I2").Value = 24.153847
$ws.Range("J2").Value = 43.5
$ws.Range("K2").Value = 24.153847
$ws.Range("L2").Value = 43.5
$ws.Range("M2").Value = 88.846153
$ws.Range("N2").Value = -269.5
$ws.Range("H111").Value = 17623.25
$ws.Range("J111").Value = 17623.25
$ws.Range("L111").Value = 17623.25
$ws.Range("N111").Value = -23757.25
$ws.Range("H126").Value = 3893.7144
$ws.Range("I126").Value = 2628
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 7884
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -5414
$ws.Range("N126").Value = -18140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 432504
$ws.Range("I14").Value = 1255002
$ws.Range("J14").Value = 21255
$ws.Range("K14").Value = 1255002
$ws.Range("L14").Value = 21255
$ws.Range("M14").Value = -1254830
$ws.Range("N14").Value = -21599
$ws.Range("H25").Value = 47072
$ws.Range("J25").Value = 47072
$ws.Range("L25").Value = 47072
$ws.Range("N25").Value = -47532
$ws.Range("H43").Value = 45208.4
$ws.Range("J43").Value = 45208.4
$ws.Range("L43").Value = 45208.4
$ws.Range("N43").Value = -45594.4
$ws.Range("H46").Value = 1663.5294
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 1828.6666
$ws.Range("K46").Value = 425
$ws.Range("L46").Value = 1828.6666
$ws.Range("M46").Value = -237
$ws.Range("N46").Value = -2204.6666
$ws.Range("H61").Value = 83337624
$ws.Range("I61").Value = 200001740
$ws.Range("J61").Value = 6114.2856
$ws.Range("K61").Value = 200001740
$ws.Range("L61").Value = 6114.2856
$ws.Range("M61").Value = -200001538
$ws.Range("N61").Value = -6518.2856
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0
$ws.Range("H109").Value = 24500
$ws.Range("J109").Value = 24500
$ws.Range("L109").Value = 24500
$ws.Range("N109").Value = -27274
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H112").Value = 67999
$ws.Range("J112").Value = 67999
$ws.Range("L112").Value = 67999
$ws.Range("N112").Value = -70953
$ws.Range("H113").Value = 83337624
$ws.Range("I113").Value = 200001740
$ws.Range("J113").Value = 6114.2856
$ws.Range("K113").Value = 200001740
$ws.Range("L113").Value = 6114.2856
$ws.Range("M113").Value = -199999570
$ws.Range("N113").Value = -10454.2856
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H116").Value = 36666.668
$ws.Range("J116").Value = 36666.668
$ws.Range("L116").Value = 36666.668
$ws.Range("N116").Value = -45844.668
$ws.Range("H118").Value = 25000
$ws.Range("J118").Value = 25000
$ws.Range("L118").Value = 25000
$ws.Range("N118").Value = -28314
$ws.Range("H119").Value = 49800
$ws.Range("J119").Value = 49800
$ws.Range("L119").Value = 49800
$ws.Range("N119").Value = -59476
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0
$ws.Range("H121").Value = 40420
$ws.Range("J121").Value = 40420
$ws.Range("L121").Value = 40420
$ws.Range("N121").Value = -43914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 7500000.5
$ws.Range("J2").Value = 6666667.5
$ws.Range("L2").Value = 6666667.5
$ws.Range("N2").Value = -6666891.5
